# Apply the bsky_followers.xlsx update:
#  - Fix three previously-missing follower counts (N12, P12, H13)
#  - Add four new tracked accounts as new columns AA:AD (header row only)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the zero/missing values that were actually real counts ---
$ws.Range("N12").Value = 8950
$ws.Range("P12").Value = 3310
$ws.Range("H13").Value = 10034

# --- Add four new accounts as additional header columns ---
$ws.Range("AA1").Value = "fuelpovertyaction.bsky.social"
$ws.Range("AB1").Value = "jrf-uk.bsky.social"
$ws.Range("AC1").Value = "e3g.bsky.social"
$ws.Range("AD1").Value = "neweconomics.bsky.social"

# Match the bold/centered header formatting already used for row 1
$ws.Range("AA1:AD1").Font.Bold = $true
$ws.Range("AA1:AD1").HorizontalAlignment = -4108

# Leave the selection where the author left it when they saved
$ws.Range("AE1").Select()
